# Closing prices instead of current
#
# The "Current Price" column (G) is being repurposed to show the previous
# closing price ("Prev Close") instead of the current/live price. Update
# the column header and the price figures, and refresh the dependent
# "Daily Change %" figures (column I) that are tied to those new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header from "Current Price" to "Prev Close"
$ws.Range("G2").Value = "Prev Close"

# Update the (previous) closing price figures in column G
$ws.Range("G3").Value = 148.98
$ws.Range("G4").Value = 139.44
$ws.Range("G6").Value = 23.6
$ws.Range("G7").Value = 398.645

# Update the daily % change text for SQ (row 6). The cell is formatted as a
# percentage, so it must be temporarily switched to text to preserve the
# literal "-0.47%" label instead of being re-interpreted as a number.
$i6 = $ws.Range("I6")
$i6.NumberFormat = "@"
$i6.Value = "-0.47%"
$i6.NumberFormat = "0%"

# Update the daily % change figure for ETH (row 7, plain numeric percentage)
$ws.Range("I7").Value = 0.1005
